$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.375.15"
$ws.Range("E2").Value = "  -3.08%  "

$ws.Range("D3").Value = "1.813.69"
$ws.Range("E3").Value = "  -3.55%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4222"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3570"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07171"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8511"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.30"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.20%  "

$ws.Range("D12").Value = "1.865.07"
$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.337"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.13%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.394"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.73%  "

$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06911"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.73"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008847"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.85%  "

$ws.Range("D21").Value = "27.831.50"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.104"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "2.125.29"
$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.969"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.72"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.08%  "

$ws.Range("E28").Value = "  -5.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.43"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.96%  "

$ws.Range("E30").Value = "  -8.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08907"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7441"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -7.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.491"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.932"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.113"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -6.09%  "

$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.078"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05231"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01912"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.764"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.92%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1647"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.20%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5015"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.334"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -8.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.266"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.21%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.35"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.62%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.28"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06424"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.95%  "

$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4615"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.608"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.36"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.12%  "

